$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update totals for Estados Unidos (row 4) ---
$ws.Range("B4").Value = 1030681
$ws.Range("C4").Value = 20325
$ws.Range("D4").Value = 140288
$ws.Range("E4").Value = 831711
$ws.Range("G4").Value = 1885
$ws.Range("H4").Value = 58682

# --- Reorder / refresh Guinea, Costa de Marfil, Republica de Yibuti, Hong Kong ---
# New order: Guinea (row86), Costa de Marfil (row87), Republica de Yibuti (row88), Hong Kong (row89)
$ws.Range("A86").Value = "Guinea"
$ws.Range("B86").Value = 1240
$ws.Range("C86").Value = 77
$ws.Range("D86").Value = 269
$ws.Range("E86").Value = 964
$ws.Range("F86").Value = 0
$ws.Range("G86").Value = 0
$ws.Range("H86").Value = 7

$ws.Range("A87").Value = "Costa de Marfil"
$ws.Range("B87").Value = 1164
$ws.Range("C87").Value = 0
$ws.Range("D87").Value = 499
$ws.Range("E87").Value = 651
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 14

$ws.Range("A88").Value = "Republica de Yibuti"
$ws.Range("B88").Value = 1072
$ws.Range("C88").Value = 37
$ws.Range("D88").Value = 498
$ws.Range("E88").Value = 572
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 2

$ws.Range("A89").Value = "Hong Kong"
$ws.Range("B89").Value = 1038
$ws.Range("C89").Value = 0
$ws.Range("D89").Value = 811
$ws.Range("E89").Value = 223
$ws.Range("F89").Value = 4
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 4

# --- Togo (row150) and Bahamas (row154) small updates ---
$ws.Range("D150").Value = 63
$ws.Range("E150").Value = 30

$ws.Range("D154").Value = 23
$ws.Range("E154").Value = 46

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 00:52"
